$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "Adjust param according to the random crop`nGrayscale`nUse border between windows"
$ws.Range("D3").Value = "Adjust param according to the random crop`nGrayscale`nUse border between windows"
$ws.Range("D4").Value = "Adjust param according to the random crop`nGrayscale`nUse border between windows`nCrop sky and shop"

$ws.Range("D5").Select()
